$d = $word.ActiveDocument

$pairs = @(
    @("71×90=6390", "26×54=1404"),
    @("39×69=2691", "38×49=1862"),
    @("97×62=6014", "56×17=952"),
    @("81×76=6156", "81×34=2754"),
    @("16×71=1136", "44×45=1980"),
    @("89×71=6319", "75×98=7350"),
    @("45×24=1080", "94×68=6392"),
    @("42×61=2562", "79×33=2607"),
    @("76×19=1444", "12×69=828"),
    @("58×23=1334", "81×76=6156"),
    @("42×51=2142", "52×29=1508"),
    @("63×96=6048", "87×32=2784"),
    @("14×72=1008", "89×79=7031"),
    @("83×71=5893", "46×56=2576"),
    @("86×37=3182", "78×66=5148"),
    @("88×46=4048", "21×95=1995"),
    @("66×61=4026", "91×41=3731"),
    @("40×16=640", "78×37=2886"),
    @("18×54=972", "39×73=2847"),
    @("34×32=1088", "74×96=7104"),
    @("93×49=4557", "92×17=1564"),
    @("52×85=4420", "40×72=2880"),
    @("73×31=2263", "97×50=4850"),
    @("27×65=1755", "99×92=9108"),
    @("75×72=5400", "42×81=3402")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
